{"js": "// Fix CI bounds in the delta (MSR) table: several cells had a stray\n// \"0.16\" placeholder in place of the real lower/upper confidence bound.\nconst replacements = [\n  { old: \"-2.86 ( 0.16 - -0.62 )\",  new: \"-2.86 ( -5.68 - -0.62 )\" },\n  { old: \"-6.01 ( 0.16 - -2.95 )\",  new: \"-6.01 ( -9.12 - -2.95 )\" },\n  { old: \"-8.87 ( 0.16 - -5.76 )\",  new: \"-8.87 ( -12.61 - -5.76 )\" },\n  { old: \"-1.56 ( -3.38 - -0.16 )\", new: \"-1.56 ( -3.45 - 0.33 )\" },\n  { old: \"4.14 ( 0.16 - 7.49 )\",    new: \"4.14 ( 0.55 - 7.49 )\" },\n  { old: \"-1.87 ( 0.16 - -0.12 )\",  new: \"-1.87 ( -3.91 - -0.12 )\" },\n  { old: \"-4.73 ( 0.16 - -3.43 )\",  new: \"-4.73 ( -6.9 - -3.43 )\" },\n  { old: \"4.79 ( 2.76 - 6.65 )\",    new: \"4.79 ( 2.86 - 6.71 )\" },\n  { old: \"6.34 ( 0.16 - 8.71 )\",    new: \"6.34 ( 4.24 - 8.71 )\" },\n  { old: \"10.48 ( 0.16 - 15.18 )\",  new: \"10.48 ( 5.81 - 15.18 )\" },\n  { old: \"4.47 ( 0.16 - 6.78 )\",    new: \"4.47 ( 2.15 - 6.78 )\" },\n  { old: \"1.61 ( 0.16 - 3.82 )\",    new: \"1.61 ( -1.19 - 3.82 )\" },\n  { old: \"5.69 ( 0.16 - 9.31 )\",    new: \"5.69 ( 2.28 - 9.31 )\" },\n  { old: \"-0.31 ( 0.16 - 1.45 )\",   new: \"-0.31 ( -1.93 - 1.45 )\" },\n  { old: \"-3.18 ( 0.16 - -1.24 )\",  new: \"-3.18 ( -5.55 - -1.24 )\" },\n];\n\nconst body = context.document.body;\n\nfor (const r of replacements) {\n  const results = body.search(r.old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(r.new, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix CI bounds in the delta (MSR) table: several cells had a stray\n# \"0.16\" placeholder in place of the real lower/upper confidence bound.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"-2.86 ( 0.16 - -0.62 )\";   New = \"-2.86 ( -5.68 - -0.62 )\" },\n    @{ Old = \"-6.01 ( 0.16 - -2.95 )\";   New = \"-6.01 ( -9.12 - -2.95 )\" },\n    @{ Old = \"-8.87 ( 0.16 - -5.76 )\";   New = \"-8.87 ( -12.61 - -5.76 )\" },\n    @{ Old = \"-1.56 ( -3.38 - -0.16 )\";  New = \"-1.56 ( -3.45 - 0.33 )\" },\n    @{ Old = \"4.14 ( 0.16 - 7.49 )\";     New = \"4.14 ( 0.55 - 7.49 )\" },\n    @{ Old = \"-1.87 ( 0.16 - -0.12 )\";   New = \"-1.87 ( -3.91 - -0.12 )\" },\n    @{ Old = \"-4.73 ( 0.16 - -3.43 )\";   New = \"-4.73 ( -6.9 - -3.43 )\" },\n    @{ Old = \"4.79 ( 2.76 - 6.65 )\";     New = \"4.79 ( 2.86 - 6.71 )\" },\n    @{ Old = \"6.34 ( 0.16 - 8.71 )\";     New = \"6.34 ( 4.24 - 8.71 )\" },\n    @{ Old = \"10.48 ( 0.16 - 15.18 )\";   New = \"10.48 ( 5.81 - 15.18 )\" },\n    @{ Old = \"4.47 ( 0.16 - 6.78 )\";     New = \"4.47 ( 2.15 - 6.78 )\" },\n    @{ Old = \"1.61 ( 0.16 - 3.82 )\";     New = \"1.61 ( -1.19 - 3.82 )\" },\n    @{ Old = \"5.69 ( 0.16 - 9.31 )\";     New = \"5.69 ( 2.28 - 9.31 )\" },\n    @{ Old = \"-0.31 ( 0.16 - 1.45 )\";    New = \"-0.31 ( -1.93 - 1.45 )\" },\n    @{ Old = \"-3.18 ( 0.16 - -1.24 )\";   New = \"-3.18 ( -5.55 - -1.24 )\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
